# Passive Skills workbook update: new skill tree entries (Farmers Guild,
# Calvary Training Grounds) plus new is_locked/is_parent columns (I, J)
# populated for all existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "is_parent" (J) values for existing rows 2-9 ---
$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("J9").Value = 0

# --- New "is_locked" (I) values for rows that did not already have one ---
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1

# --- Row 10: Farmers Guild ---
$ws.Range("A10").Value = "Farmers Guild"
$ws.Range("B10").Value = "At level 5, and not before, you will get, at hourly reset, 100% of your maximum population as you have learned how to use food efficiently to feed your population."
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 4
$ws.Range("G10").Value = "Building Research"
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 0

# --- Row 11: Calvary Training Grounds ---
$ws.Range("A11").Value = "Calvary Training Grounds"
$ws.Range("B11").Value = "This will unlock the Calvary Training Grounds Building, allowing you to recruit Mounted Knights and Mounted Archers."
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 10
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = "Farmers Guild"
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0

# --- Selection moved to J3 as part of the re-save ---
$ws.Range("J3").Select()
